# Compatibilização com novo formato da base da ANEEL
# Adds a new "classe" (Consumo Próprio) with its AT/BT, local/remoto
# combinations to the segmento lookup table, and extends the
# AutoFilter/defined-name ranges accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the four new data rows (30-33) ---------------------------------
$ws.Range("A30").Value = "Consumo Próprio"
$ws.Range("B30").Value = "AT"
$ws.Range("C30").Value = "local"
$ws.Range("D30").Value = "comercial_at"

$ws.Range("A31").Value = "Consumo Próprio"
$ws.Range("B31").Value = "AT"
$ws.Range("C31").Value = "remoto"
$ws.Range("D31").Value = "comercial_at_remoto"

$ws.Range("A32").Value = "Consumo Próprio"
$ws.Range("B32").Value = "BT"
$ws.Range("C32").Value = "local"
$ws.Range("D32").Value = "comercial_bt"

$ws.Range("A33").Value = "Consumo Próprio"
$ws.Range("B33").Value = "BT"
$ws.Range("C33").Value = "remoto"
$ws.Range("D33").Value = "comercial_bt"

# --- Widen column A (no longer auto "best fit") -------------------------
$ws.Columns.Item(1).ColumnWidth = 16.3

# --- Re-apply the AutoFilter over the new, larger range -----------------
$ws.AutoFilterMode = $false
$ws.Range("A1:D33").AutoFilter() | Out-Null

# Keep the _FilterDatabase defined name in sync with the new range
$wb.Names("_xlnm._FilterDatabase").Value = "=Sheet1!`$A`$1:`$D`$33"

# --- Match the cursor/selection left behind by the editor ---------------
$ws.Range("A30").Select() | Out-Null
